$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Row 10: "sub menu silder bar: https://youtu.be/E5NB9crfQZs?si=Q8c77KPhtsqENjsh "
$text10 = "sub menu silder bar:" + " https://youtu.be/E5NB9crfQZs?si=Q8c77KPhtsqENjsh "
$ws.Range("B10").Value = $text10
$run10 = $ws.Range("B10").Characters(21, 50)
$run10.Font.Bold = $true

# Row 11: "main : https://www.youtube.com/watch?v=CkHyDYeImjY"
$text11 = "main :" + " https://www.youtube.com/watch?v=CkHyDYeImjY"
$ws.Range("B11").Value = $text11
$run11 = $ws.Range("B11").Characters(7, 44)
$run11.Font.Bold = $true

# Update selection to match the edited sheet view (B11 selected)
[void]$ws.Range("B11").Select()
